$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff (cryptos list refresh, Thu May 30 11:11:58 UTC 2024).
# A leading apostrophe forces Excel to store the value as literal text (quote-prefixed)
# instead of re-typing numeric-looking strings (e.g. "594.34", "0.0000263") as numbers.
$ws.Range('D2').Value = ("'" + '67.921.44')
$ws.Range('E2').Value = ("'" + '  +0.10%  ')
$ws.Range('D3').Value = ("'" + '3.745.23')
$ws.Range('E3').Value = ("'" + '  -1.91%  ')
$ws.Range('E4').Value = ("'" + '  +0.01%  ')
$ws.Range('D5').Value = ("'" + '594.34')
$ws.Range('E5').Value = ("'" + '  -0.91%  ')
$ws.Range('D6').Value = ("'" + '166.52')
$ws.Range('E6').Value = ("'" + '  -1.72%  ')
$ws.Range('D7').Value = ("'" + '3.742.77')
$ws.Range('E7').Value = ("'" + '  -1.98%  ')
$ws.Range('E8').Value = ("'" + '  +0.02%  ')
$ws.Range('D9').Value = ("'" + '0.521')
$ws.Range('E9').Value = ("'" + '  -1.66%  ')
$ws.Range('E10').Value = ("'" + '  -3.51%  ')
$ws.Range('D11').Value = ("'" + '6.49')
$ws.Range('E11').Value = ("'" + '  -0.66%  ')
$ws.Range('D12').Value = ("'" + '0.451')
$ws.Range('E12').Value = ("'" + '  -2.48%  ')
$ws.Range('D13').Value = ("'" + '0.0000263')
$ws.Range('E13').Value = ("'" + '  -4.23%  ')
$ws.Range('D14').Value = ("'" + '36.41')
$ws.Range('E14').Value = ("'" + '  -1.64%  ')
$ws.Range('D15').Value = ("'" + '4.371.63')
$ws.Range('E15').Value = ("'" + '  -2.11%  ')
$ws.Range('D16').Value = ("'" + '3.738.77')
$ws.Range('E16').Value = ("'" + '  -2.28%  ')
$ws.Range('D17').Value = ("'" + '67.903.78')
$ws.Range('E17').Value = ("'" + '  -0.06%  ')
$ws.Range('D18').Value = ("'" + '18.26')
$ws.Range('E18').Value = ("'" + '  -1.96%  ')
$ws.Range('E19').Value = ("'" + '  -5.51%  ')
$ws.Range('E20').Value = ("'" + '  -0.61%  ')
$ws.Range('D21').Value = ("'" + '10.78')
$ws.Range('E21').Value = ("'" + '  -0.80%  ')
$ws.Range('D22').Value = ("'" + '468.60')
$ws.Range('E22').Value = ("'" + '  -0.07%  ')
$ws.Range('D23').Value = ("'" + '0.703')
$ws.Range('E23').Value = ("'" + '  -5.46%  ')
$ws.Range('D24').Value = ("'" + '83.05')
$ws.Range('E24').Value = ("'" + '  -0.73%  ')
$ws.Range('B25').Value = ("'" + 'PEPE')
$ws.Range('C25').Value = ("'" + 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe')
$ws.Range('D25').Value = ("'" + '0.0000136')
$ws.Range('E25').Value = ("'" + '  -10.35%  ')
$ws.Range('B26').Value = ("'" + 'Fetch.AI')
$ws.Range('C26').Value = ("'" + 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet')
$ws.Range('D26').Value = ("'" + '2.22')
$ws.Range('E26').Value = ("'" + '  -3.79%  ')
$ws.Range('E27').Value = ("'" + '  -1.70%  ')
$ws.Range('D28').Value = ("'" + '10.16')
$ws.Range('E28').Value = ("'" + '  -2.08%  ')
$ws.Range('E29').Value = ("'" + '  -0.02%  ')
$ws.Range('D30').Value = ("'" + '3.890.08')
$ws.Range('E30').Value = ("'" + '  -1.96%  ')
$ws.Range('D31').Value = ("'" + '2.79')
$ws.Range('E31').Value = ("'" + '  -4.86%  ')
$ws.Range('E32').Value = ("'" + '  -1.99%  ')
$ws.Range('E33').Value = ("'" + '  -5.41%  ')
$ws.Range('D34').Value = ("'" + '29.84')
$ws.Range('E34').Value = ("'" + '  -3.06%  ')
$ws.Range('D35').Value = ("'" + '0.999')
$ws.Range('D36').Value = ("'" + '9.07')
$ws.Range('E36').Value = ("'" + '  -2.99%  ')
$ws.Range('D37').Value = ("'" + '3.696.49')
$ws.Range('E37').Value = ("'" + '  -2.34%  ')
$ws.Range('E38').Value = ("'" + '  -4.60%  ')
$ws.Range('D39').Value = ("'" + '3.47')
$ws.Range('E39').Value = ("'" + '  -9.81%  ')
$ws.Range('E40').Value = ("'" + '  -1.13%  ')
$ws.Range('D41').Value = ("'" + '0.991')
$ws.Range('E41').Value = ("'" + '  -2.03%  ')
$ws.Range('E42').Value = ("'" + '  -3.93%  ')
$ws.Range('E43').Value = ("'" + '  -0.09%  ')
$ws.Range('E44').Value = ("'" + '  +0.02%  ')
$ws.Range('E45').Value = ("'" + '  -3.59%  ')
$ws.Range('D46').Value = ("'" + '8.60')
$ws.Range('E46').Value = ("'" + '  -2.28%  ')
$ws.Range('E47').Value = ("'" + '  -2.79%  ')
$ws.Range('D48').Value = ("'" + '45.35')
$ws.Range('E48').Value = ("'" + '  -2.68%  ')
$ws.Range('D49').Value = ("'" + '393.35')
$ws.Range('E49').Value = ("'" + '  -4.26%  ')
$ws.Range('D50').Value = ("'" + '143.91')
$ws.Range('E50').Value = ("'" + '  +0.80%  ')
$ws.Range('D51').Value = ("'" + '25.45')
$ws.Range('E51').Value = ("'" + '  -0.23%  ')
